$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "B2" = 0.2291666666666667
    "C2" = 0.5083333333333333
    "J2" = 0.01666666666666667
    "P2" = 0.175
    "S2" = 0.07083333333333333
    "B3" = 0.01652892561983471
    "C3" = 0.02479338842975207
    "J3" = 0.03305785123966942
    "P3" = 0.7768595041322314
    "S3" = 0.1487603305785124
    "J4" = 0.02631578947368421
    "P4" = 0.7631578947368421
    "S4" = 0.2105263157894737
    "P5" = 0.75
    "S5" = 0.25
    "B6" = 0.05376344086021505
    "D6" = 0.03225806451612903
    "F6" = 0.03225806451612903
    "J6" = 0.1720430107526882
    "O6" = 0.005376344086021506
    "Q6" = 0.1989247311827957
    "R6" = 0.08064516129032258
    "S6" = 0.4247311827956989
    "B7" = 0.09523809523809523
    "D7" = 0.02380952380952381
    "E7" = 0.004761904761904762
    "F7" = 0.0380952380952381
    "J7" = 0.1476190476190476
    "O7" = 0.01904761904761905
    "Q7" = 0.1333333333333333
    "R7" = 0.1238095238095238
    "S7" = 0.4142857142857143
    "B8" = 0.07713498622589532
    "D8" = 0.01377410468319559
    "F8" = 0.05785123966942149
    "J8" = 0.140495867768595
    "O8" = 0.01377410468319559
    "Q8" = 0.1570247933884298
    "R8" = 0.09917355371900827
    "S8" = 0.440771349862259
    "B9" = 0.1156462585034014
    "D9" = 0.0272108843537415
    "E9" = 0.006802721088435374
    "F9" = 0.06802721088435375
    "J9" = 0.09523809523809523
    "Q9" = 0.1292517006802721
    "R9" = 0.08163265306122448
    "S9" = 0.4761904761904762
    "B10" = 0.1058091286307054
    "D10" = 0.01867219917012448
    "E10" = 0.002074688796680498
    "F10" = 0.0954356846473029
    "J10" = 0.1047717842323651
    "O10" = 0.01867219917012448
    "Q10" = 0.1670124481327801
    "R10" = 0.09854771784232365
    "S10" = 0.3890041493775934
    "G11" = 0.1467065868263473
    "J11" = 0.07784431137724551
    "K11" = 0.2305389221556886
    "L11" = 0.5209580838323353
    "S11" = 0.02395209580838323
    "G12" = 0.7457627118644068
    "J12" = 0.1977401129943503
    "L12" = 0.01694915254237288
    "S12" = 0.03954802259887006
    "F13" = 0.02083333333333333
    "G13" = 0.6875
    "J13" = 0.1458333333333333
    "S13" = 0.1458333333333333
    "F15" = 0.02339181286549707
    "H15" = 0.1695906432748538
    "I15" = 0.03508771929824561
    "J15" = 0.3333333333333333
    "K15" = 0.07602339181286549
    "M15" = 0.01754385964912281
    "O15" = 0.07602339181286549
    "S15" = 0.2690058479532164
    "F16" = 0.006289308176100629
    "H16" = 0.1383647798742138
    "I16" = 0.09433962264150944
    "J16" = 0.3710691823899371
    "K16" = 0.1069182389937107
    "M16" = 0.02515723270440252
    "O16" = 0.06289308176100629
    "S16" = 0.1949685534591195
    "F17" = 0.006600660066006601
    "H17" = 0.1749174917491749
    "I17" = 0.0594059405940594
    "J17" = 0.3927392739273927
    "K17" = 0.1188118811881188
    "M17" = 0.0264026402640264
    "O17" = 0.03630363036303631
    "S17" = 0.1848184818481848
    "F18" = 0.005376344086021506
    "H18" = 0.1451612903225807
    "I18" = 0.09139784946236559
    "J18" = 0.2688172043010753
    "K18" = 0.1344086021505376
    "M18" = 0.02150537634408602
    "N18" = 0.005376344086021506
    "O18" = 0.06989247311827956
    "S18" = 0.2580645161290323
    "F19" = 0.01221498371335505
    "H19" = 0.1872964169381107
    "I19" = 0.07573289902280131
    "J19" = 0.3200325732899023
    "K19" = 0.1327361563517915
    "M19" = 0.02361563517915309
    "O19" = 0.06270358306188925
    "S19" = 0.1856677524429967
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
